$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $values) {
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}

# Row 2: ECs -> Thpo -> Mpl -> ECs (values updated)
Set-Row 2 @("ECs", "Thpo", "Mpl", "ECs", 1, 0.3333333333333333, 0.057981, 0.173943, 0.2608001391389875, 0.2608001391389875, 3, 1, 0.7852063333333333, 2.355619, 0.984218618779221, 0.984218618779221, 0.045527048413, 0.4097434357169999, 0.256684352720803, 0.256684352720803)

# Row 3: now ECs -> Thpo -> Mpl -> Resolving-Mac
Set-Row 3 @("ECs", "Thpo", "Mpl", "Resolving-Mac", 1, 0.3333333333333333, 0.057981, 0.173943, 0.2608001391389875, 0.2608001391389875, 1, 0.3333333333333333, 0.01259033333333333, 0.037771, 0.0157813812207789, 0.01578138122077889, 0.000730000117, 0.006570001052999999, 0.004115786418184541, 0.00411578641818454)

# Row 4: now MuSCs -> Thpo -> Mpl -> ECs
Set-Row 4 @("MuSCs", "Thpo", "Mpl", "ECs", 3, 1, 0.1643386666666667, 0.493016, 0.7391998608610124, 0.7391998608610125, 3, 1, 0.7852063333333333, 2.355619, 0.984218618779221, 0.984218618779221, 0.1290397618782222, 1.161357856904, 0.727534266058418, 0.7275342660584181)

# Row 5 (new): MuSCs -> Thpo -> Mpl -> Resolving-Mac
Set-Row 5 @("MuSCs", "Thpo", "Mpl", "Resolving-Mac", 3, 1, 0.1643386666666667, 0.493016, 0.7391998608610124, 0.7391998608610125, 1, 0.3333333333333333, 0.01259033333333333, 0.037771, 0.0157813812207789, 0.01578138122077889, 0.002069078592888889, 0.018621707336, 0.01166559480259435, 0.01166559480259435)
